$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values like "1.00" or "58.409.82" that Excel
# would otherwise coerce into numbers (dropping trailing zeros, collapsing
# multi-dot "thousands" groupings, etc). Format it as Text first so the
# values are stored verbatim, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '58.409.82'
$ws.Range("E2").Value = '  +0.90%  '

# Row 3
$ws.Range("D3").Value = '3.147.01'
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '536.22'
$ws.Range("E5").Value = '  +1.10%  '

# Row 6
$ws.Range("D6").Value = '139.85'
$ws.Range("E6").Value = '  +0.83%  '

# Row 7
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").Value = '3.144.94'
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  +4.69%  '

# Row 10
$ws.Range("D10").Value = '7.31'
$ws.Range("E10").Value = '  +1.21%  '

# Row 11
$ws.Range("D11").Value = '0.108'
$ws.Range("E11").Value = '  -0.22%  '

# Row 12
$ws.Range("D12").Value = '0.416'
$ws.Range("E12").Value = '  +4.17%  '

# Row 13
$ws.Range("D13").Value = '3.663.07'
$ws.Range("E13").Value = '  -0.54%  '

# Row 14
$ws.Range("E14").Value = '  +1.35%  '

# Row 15
$ws.Range("D15").Value = '25.74'
$ws.Range("E15").Value = '  +0.53%  '

# Row 16
$ws.Range("D16").Value = '0.0000166'
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("D17").Value = '58.469.15'
$ws.Range("E17").Value = '  +0.79%  '

# Row 18
$ws.Range("D18").Value = '3.132.28'
$ws.Range("E18").Value = '  -0.28%  '

# Row 19
$ws.Range("D19").Value = '6.08'
$ws.Range("E19").Value = '  +1.25%  '

# Row 20
$ws.Range("D20").Value = '12.75'
$ws.Range("E20").Value = '  -0.17%  '

# Row 21
$ws.Range("D21").Value = '8.20'
$ws.Range("E21").Value = '  +2.72%  '

# Row 22
$ws.Range("D22").Value = '361.18'
$ws.Range("E22").Value = '  +2.22%  '

# Row 23
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("D24").Value = '69.28'
$ws.Range("E24").Value = '  +0.78%  '

# Row 25
$ws.Range("D25").Value = '0.508'
$ws.Range("E25").Value = '  -0.08%  '

# Row 26
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -1.24%  '

# Row 27
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.30%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0884'
$ws.Range("E28").Value = '  -3.95%  '

# Row 29
$ws.Range("D29").Value = '7.36'
$ws.Range("E29").Value = '  -2.69%  '

# Row 30
$ws.Range("D30").Value = '6.18'
$ws.Range("E30").Value = '  +0.00%  '

# Row 31
$ws.Range("D31").Value = '1.89'
$ws.Range("E31").Value = '  -0.29%  '

# Row 32
$ws.Range("D32").Value = '21.73'
$ws.Range("E32").Value = '  +2.42%  '

# Row 33
$ws.Range("D33").Value = '5.15'
$ws.Range("E33").Value = '  +2.32%  '

# Row 34
$ws.Range("D34").Value = '1.16'
$ws.Range("E34").Value = '  -2.48%  '

# Row 35
$ws.Range("D35").Value = '159.25'
$ws.Range("E35").Value = '  +0.84%  '

# Row 36
$ws.Range("D36").Value = '6.12'
$ws.Range("E36").Value = '  -1.20%  '

# Row 37
$ws.Range("D37").Value = '25.93'
$ws.Range("E37").Value = '  -2.30%  '

# Row 38
$ws.Range("D38").Value = '1.28'
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("D39").Value = '1.70'
$ws.Range("E39").Value = '  +4.16%  '

# Row 40
$ws.Range("D40").Value = '0.0675'
$ws.Range("E40").Value = '  +0.43%  '

# Row 41
$ws.Range("D41").Value = '2.502.65'
$ws.Range("E41").Value = '  +6.80%  '

# Row 42
$ws.Range("D42").Value = '0.704'
$ws.Range("E42").Value = '  -0.32%  '

# Row 43
$ws.Range("D43").Value = '4.03'
$ws.Range("E43").Value = '  -4.28%  '

# Row 44
$ws.Range("D44").Value = '37.61'
$ws.Range("E44").Value = '  +2.62%  '

# Row 45
$ws.Range("D45").Value = '3.180.49'
$ws.Range("E45").Value = '  -0.16%  '

# Row 46
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0269'
$ws.Range("E46").Value = '  -0.73%  '

# Row 47
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +3.16%  '

# Row 49
$ws.Range("D49").Value = '6.10'
$ws.Range("E49").Value = '  +0.77%  '

# Row 50
$ws.Range("D50").Value = '20.00'
$ws.Range("E50").Value = '  -2.07%  '

# Row 51
$ws.Range("D51").Value = '0.749'
$ws.Range("E51").Value = '  -2.70%  '
